# High-Profile Narrow Fillister Head Phillips Screws
# Insert a new header row of column indices (0-11) above the existing
# text header row, pushing the text header row (and all data rows)
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing header row (row 1) text before we touch anything.
$headers = @()
for ($c = 1; $c -le 12; $c++) {
    $headers += $ws.Cells.Item(1, $c).Value2
}

# Insert a brand-new row at position 2: this duplicates row 1's
# formatting into the new row 2 and shifts the old rows 2..37 down to
# 3..38, while leaving row 1 (and its style) untouched.
$ws.Rows.Item(2).Insert()

# The newly inserted row 2 picked up row 1's bold/border/centered
# style; the target file has it as a plain, unstyled row (matching the
# rest of the data rows), so strip that back off.
$ws.Range("A2:L2").ClearFormats()

# Move the old header text down into row 2 (I2 stays blank, same as
# the old I1; K2/L2 are left blank too).
$ws.Cells.Item(2, 1).Value = $headers[0]
$ws.Cells.Item(2, 2).Value = $headers[1]
$ws.Cells.Item(2, 3).Value = $headers[2]
$ws.Cells.Item(2, 4).Value = $headers[3]
$ws.Cells.Item(2, 5).Value = $headers[4]
$ws.Cells.Item(2, 6).Value = $headers[5]
$ws.Cells.Item(2, 7).Value = $headers[6]
$ws.Cells.Item(2, 8).Value = $headers[7]
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 10).Value = $headers[9]
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""

# Row 1 becomes a purely numeric index row (0-11), keeping its
# original bold/border/centered style (s="1") since we never cleared
# it.
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}
